$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row at 107, shifting existing rows 107+ down by one.
$ws.Rows.Item(107).Insert()

# Populate the new row 107 with data. The "category" columns (A-C, E-J, R)
# are copied from the entry that used to occupy row 107 (now shifted to row 108).
$ws.Range("A107").Value2 = $ws.Range("A108").Value2
$ws.Range("B107").Value2 = $ws.Range("B108").Value2
$ws.Range("C107").Value2 = $ws.Range("C108").Value2
$ws.Range("D107").Value2 = 44484
$ws.Range("E107").Value2 = $ws.Range("E108").Value2
$ws.Range("F107").Value2 = $ws.Range("F108").Value2
$ws.Range("G107").Value2 = $ws.Range("G108").Value2
$ws.Range("H107").Value2 = $ws.Range("H108").Value2
$ws.Range("I107").Value2 = $ws.Range("I108").Value2
$ws.Range("J107").Value2 = $ws.Range("J108").Value2
$ws.Range("K107").Value2 = 15000
$ws.Range("L107").Value2 = 15000
$ws.Range("M107").Value2 = 15000
$ws.Range("N107").Value2 = "$/caja 60 unidades"
$ws.Range("O107").Value2 = "Región de O'Higgins"
$ws.Range("P107").Value2 = 250
$ws.Range("Q107").Value2 = 60
$ws.Range("R107").Value2 = $ws.Range("R108").Value2
